# Applies the WorkAbsencePermitType rename + supporting text/content edits.

$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1) Entity name rename: BusinessTripCostComponent -> WorkAbsencePermitType
#    (occurs 9x across the document: title, table rows, code samples)
Replace-All "BusinessTripCostComponent" "WorkAbsencePermitType"

# 2) Description sentence: "Menghapus Data Komponen Biaya Perjalanan Bisnis"
#    -> "Menghapus Data Izin Jenis Absensi Ketidakhadiran"
Replace-All "Menghapus Data Komponen Biaya Perjalanan Bisnis" "Menghapus Data Izin Jenis Absensi Ketidakhadiran"

# 3) Date fix: November 23th -> November 25th
Replace-All "23" "25"

# 4) Refreshed sample JWT bearer token
Replace-All "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJsb2dnZWRJbkFzIjoidGVndWgucHJhdGFtYSIsImlhdCI6MTYwNjA5Nzg4MH0.d1AB_XF31WOFS7dhxvEHyJmPybR5ju4YHiuF_ZbSf5Q" `
             "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJsb2dnZWRJbkFzIjoidGVndWgucHJhdGFtYSIsImlhdCI6MTYwNjI2OTA1NH0.NjJJegg6WRVQ3LHksbKcni92MkyzjfYpxzrFvgLu2FQ"

# 5) Sample recordID: 810000000000001 -> 410000000000001
Replace-All "81" "41"
